# ----------------------------------------------------------------------
# Applies the three changes described by the diff:
#  1. Merge the two runs that used to be split by the "_GoBack" bookmark
#     in the judges (评委) paragraph back into a single run (the bookmark
#     no longer lives there).
#  2. Insert the character "不" in the "一出门..." paragraph, right
#     before "一开始", and relocate the "_GoBack" bookmark to that new
#     edit point (splitting the run there, just like Word does whenever
#     the cursor was last left in the middle of a run).
#  3. Mark the built-in "Normal Table" style as a Quick Style
#     (w:qFormat), matching the QuickStyle flag toggled in styles.xml.
# ----------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1) Merge "：李春华（A）、殷文斐（" + "B）、姜烨（C）、..." ------------
# Removing the bookmark alone does not coalesce the two runs, so we
# delete the bookmark, then delete the last character right before its
# old location and retype it; Word then re-emits the whole contiguous,
# identically formatted span as a single run.
$oldBm = $d.Bookmarks("_GoBack")
$oldBmStart = $oldBm.Start
$oldBm.Delete()

$lastCharBeforeBm = $d.Range($oldBmStart - 1, $oldBmStart)
$lastCharText = $lastCharBeforeBm.Text
$lastCharBeforeBm.Delete()
$lastCharBeforeBm.InsertAfter($lastCharText)

# --- 2) Insert "不" before "一开始" and move "_GoBack" there -------------
$searchRange = $d.Content.Duplicate
$searchRange.Find.Execute("为什么我一开始就作答辩主讲")

$editPos = $searchRange.Start + 4
$editPoint = $d.Range($editPos, $editPos)
$editPoint.InsertAfter("不")

# Re-adding a bookmark named "_GoBack" automatically relocates it (bookmark
# names are unique), so the stale one from step 1's paragraph disappears.
$newBmPos = $editPos + 1
$newBmRange = $d.Range($newBmPos, $newBmPos)
$d.Bookmarks.Add("_GoBack", $newBmRange)

# --- 3) Flag "Normal Table" as a Quick Style -----------------------------
$tableStyle = $d.Styles("Normal Table")
$tableStyle.QuickStyle = $true
